$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header column (H1), copying the format used by the
# existing header cells (e.g. G1: bold, centered, bordered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the new column's data values.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
